$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("Tableau1")
$col = $lo.ListColumns.Item(3)
$col.Foo()
